$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cost ($)" column (L) for the meter usage rows can now be reported
# as "Not Available" (text) instead of a numeric cost value, matching how
# Portfolio Manager meter imports surface this field.
$range = $ws.Range("L7:L14")
$range.Value = "Not Available"
$range.NumberFormat = "@"
$range.HorizontalAlignment = 1
